$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.933
$ws.Range("A21").Value = -20.682
$ws.Range("A23").Value = -20.808
$ws.Range("B24").Value = 5.218999999999999
$ws.Range("A25").Value = -21.589
$ws.Range("B28").Value = 5.468000000000001
$ws.Range("B36").Value = 7.363
$ws.Range("B45").Value = 5.832
$ws.Range("B48").Value = 5.609
$ws.Range("B49").Value = 6.237
$ws.Range("B52").Value = 5.968000000000001
$ws.Range("A53").Value = -20.335
$ws.Range("B53").Value = 8.151
$ws.Range("B54").Value = 5.197
$ws.Range("A57").Value = -22.178
$ws.Range("A59").Value = -22.461
$ws.Range("A69").Value = -21.476
$ws.Range("B70").Value = 4.99
$ws.Range("A79").Value = -21.305
$ws.Range("A83").Value = -21.976
$ws.Range("B86").Value = 5.197
$ws.Range("B87").Value = 4.653
$ws.Range("A93").Value = -21.533
$ws.Range("B101").Value = 5.217000000000001
